$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 43.91073733333334
$ws.Range("H2").Value = 131.732212
$ws.Range("I2").Value = 0.4010337406460291
$ws.Range("J2").Value = 0.4010337406460291
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 22.529461
$ws.Range("N2").Value = 67.58838299999999
$ws.Range("O2").Value = 0.3889626001872417
$ws.Range("P2").Value = 0.3889626001872417
$ws.Range("Q2").Value = 989.2852442325773
$ws.Range("R2").Value = 8903.567198093195
$ws.Range("S2").Value = 0.1559871265244954
$ws.Range("T2").Value = 0.1559871265244954
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 43.91073733333334
$ws.Range("H3").Value = 131.732212
$ws.Range("I3").Value = 0.4010337406460291
$ws.Range("J3").Value = 0.4010337406460291
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 23.19370566666667
$ws.Range("N3").Value = 69.58111700000001
$ws.Range("O3").Value = 0.4004305324518962
$ws.Range("P3").Value = 0.4004305324518962
$ws.Range("Q3").Value = 1018.452717315645
$ws.Range("R3").Value = 9166.074455840806
$ws.Range("S3").Value = 0.1605861542980651
$ws.Range("T3").Value = 0.1605861542980651
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 43.91073733333334
$ws.Range("H4").Value = 131.732212
$ws.Range("I4").Value = 0.4010337406460291
$ws.Range("J4").Value = 0.4010337406460291
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 12.19875433333333
$ws.Range("N4").Value = 36.596263
$ws.Range("O4").Value = 0.2106068673608621
$ws.Range("P4").Value = 0.2106068673608621
$ws.Range("Q4").Value = 535.6562973248618
$ws.Range("R4").Value = 4820.906675923756
$ws.Range("S4").Value = 0.08446045982346861
$ws.Range("T4").Value = 0.08446045982346863
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 30.65522666666666
$ws.Range("H5").Value = 91.96567999999999
$ws.Range("I5").Value = 0.2799720744190927
$ws.Range("J5").Value = 0.2799720744190927
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 22.529461
$ws.Range("N5").Value = 67.58838299999999
$ws.Range("O5").Value = 0.3889626001872417
$ws.Range("P5").Value = 0.3889626001872417
$ws.Range("Q5").Value = 690.6457336328266
$ws.Range("R5").Value = 6215.811602695439
$ws.Range("S5").Value = 0.1088986660458662
$ws.Range("T5").Value = 0.1088986660458662
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 30.65522666666666
$ws.Range("H6").Value = 91.96567999999999
$ws.Range("I6").Value = 0.2799720744190927
$ws.Range("J6").Value = 0.2799720744190927
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 23.19370566666667
$ws.Range("N6").Value = 69.58111700000001
$ws.Range("O6").Value = 0.4004305324518962
$ws.Range("P6").Value = 0.4004305324518962
$ws.Range("Q6").Value = 711.0083044516178
$ws.Range("R6").Value = 6399.07474006456
$ws.Range("S6").Value = 0.1121093668312992
$ws.Range("T6").Value = 0.1121093668312992
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 30.65522666666666
$ws.Range("H7").Value = 91.96567999999999
$ws.Range("I7").Value = 0.2799720744190927
$ws.Range("J7").Value = 0.2799720744190927
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 12.19875433333333
$ws.Range("N7").Value = 36.596263
$ws.Range("O7").Value = 0.2106068673608621
$ws.Range("P7").Value = 0.2106068673608621
$ws.Range("Q7").Value = 373.9555791393155
$ws.Range("R7").Value = 3365.60021225384
$ws.Range("S7").Value = 0.05896404154192727
$ws.Range("T7").Value = 0.05896404154192728
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 34.92790866666667
$ws.Range("H8").Value = 104.783726
$ws.Range("I8").Value = 0.3189941849348781
$ws.Range("J8").Value = 0.3189941849348781
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 22.529461
$ws.Range("N8").Value = 67.58838299999999
$ws.Range("O8").Value = 0.3889626001872417
$ws.Range("P8").Value = 0.3889626001872417
$ws.Range("Q8").Value = 786.9069561172286
$ws.Range("R8").Value = 7082.162605055058
$ws.Range("S8").Value = 0.12407680761688
$ws.Range("T8").Value = 0.12407680761688
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 34.92790866666667
$ws.Range("H9").Value = 104.783726
$ws.Range("I9").Value = 0.3189941849348781
$ws.Range("J9").Value = 0.3189941849348781
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 23.19370566666667
$ws.Range("N9").Value = 69.58111700000001
$ws.Range("O9").Value = 0.4004305324518962
$ws.Range("P9").Value = 0.4004305324518962
$ws.Range("Q9").Value = 810.1076331668826
$ws.Range("R9").Value = 7290.968698501943
$ws.Range("S9").Value = 0.1277350113225319
$ws.Range("T9").Value = 0.1277350113225319
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 34.92790866666667
$ws.Range("H10").Value = 104.783726
$ws.Range("I10").Value = 0.3189941849348781
$ws.Range("J10").Value = 0.3189941849348781
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 12.19875433333333
$ws.Range("N10").Value = 36.596263
$ws.Range("O10").Value = 0.2106068673608621
$ws.Range("P10").Value = 0.2106068673608621
$ws.Range("Q10").Value = 426.0769772017709
$ws.Range("R10").Value = 3834.692794815938
$ws.Range("S10").Value = 0.06718236599546619
$ws.Range("T10").Value = 0.0671823659954662
